$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 18.46262876487648
$ws.Range("B3").Value  = 0.9993073109016987
$ws.Range("B4").Value  = 13105.86444722955
$ws.Range("B5").Value  = 11888.26840371751
$ws.Range("B6").Value  = 266.2820652073141
$ws.Range("B7").Value  = 240.280632614111
$ws.Range("B8").Value  = 20217.60000000001
$ws.Range("B9").Value  = 4145.371415078048
$ws.Range("B10").Value = 213199.8946942892
$ws.Range("B11").Value = 0.1067881886010138
$ws.Range("B12").Value = 0.3110726595435515
$ws.Range("B13").Value = 0.3500000000000072
$ws.Range("B14").Value = 0.997812849465509
$ws.Range("B15").Value = 0.4700829277784451
